$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches original inlineStr/text cell type) before assigning
# numeric-looking literal strings, so Excel does not auto-convert them to numbers.
$ws.Range("D2:D50").NumberFormat = "@"
$ws.Range("E2:E50").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "300.71"
$ws.Range("G2").Value = "9"
$ws.Range("D3").Value = "31.47"
$ws.Range("E3").Value = "0.64%"
$ws.Range("G3").Value = "9"
$ws.Range("D4").Value = "5.086"
$ws.Range("E4").Value = "-1.28%"
$ws.Range("G4").Value = "9"
$ws.Range("D5").Value = "0.07851"
$ws.Range("E5").Value = "-3.12%"
$ws.Range("G5").Value = "9"
$ws.Range("D6").Value = "2.345"
$ws.Range("E6").Value = "-6.07%"
$ws.Range("G6").Value = "9"
$ws.Range("D7").Value = "7.813"
$ws.Range("E7").Value = "-0.13%"
$ws.Range("G7").Value = "9"
$ws.Range("D8").Value = "3.824"
$ws.Range("E8").Value = "-0.71%"
$ws.Range("G8").Value = "9"
$ws.Range("D9").Value = "0.9180"
$ws.Range("E9").Value = "0.59%"
$ws.Range("G9").Value = "9"
$ws.Range("D10").Value = "0.1746"
$ws.Range("E10").Value = "2.02%"
$ws.Range("G10").Value = "9"
$ws.Range("D11").Value = "0.07576"
$ws.Range("E11").Value = "2.84%"
$ws.Range("G11").Value = "9"
$ws.Range("D12").Value = "0.09293"
$ws.Range("E12").Value = "16.37%"
$ws.Range("G12").Value = "9"
$ws.Range("D13").Value = "0.02990"
$ws.Range("E13").Value = "-1.67%"
$ws.Range("G13").Value = "9"
$ws.Range("D14").Value = "0.1003"
$ws.Range("E14").Value = "0.73%"
$ws.Range("G14").Value = "9"
$ws.Range("E15").Value = "0.63%"
$ws.Range("G15").Value = "9"
$ws.Range("D16").Value = "0.005941"
$ws.Range("E16").Value = "-1.84%"
$ws.Range("G16").Value = "9"
$ws.Range("D17").Value = "3.472"
$ws.Range("E17").Value = "-0.69%"
$ws.Range("G17").Value = "9"
$ws.Range("E18").Value = "0.29%"
$ws.Range("G18").Value = "9"
$ws.Range("E19").Value = "-0.92%"
$ws.Range("G19").Value = "9"
$ws.Range("E20").Value = "-4.64%"
$ws.Range("G20").Value = "9"
$ws.Range("D21").Value = "4.029"
$ws.Range("E21").Value = "-12.82%"
$ws.Range("G21").Value = "9"
$ws.Range("D22").Value = "0.1699"
$ws.Range("E22").Value = "5.65%"
$ws.Range("G22").Value = "9"
$ws.Range("D23").Value = "0.04615"
$ws.Range("E23").Value = "0.26%"
$ws.Range("G23").Value = "9"
$ws.Range("D24").Value = "0.001247"
$ws.Range("E24").Value = "-1.56%"
$ws.Range("G24").Value = "9"
$ws.Range("E25").Value = "0.49%"
$ws.Range("G25").Value = "9"
$ws.Range("D26").Value = "0.0001248"
$ws.Range("E26").Value = "5.36%"
$ws.Range("G26").Value = "9"
$ws.Range("E27").Value = "-1.87%"
$ws.Range("G27").Value = "9"
$ws.Range("G28").Value = "9"
$ws.Range("G29").Value = "9"
$ws.Range("G30").Value = "9"
$ws.Range("G31").Value = "9"
$ws.Range("G32").Value = "9"
$ws.Range("G33").Value = "9"
$ws.Range("G34").Value = "9"
$ws.Range("G35").Value = "9"
$ws.Range("G36").Value = "9"
$ws.Range("G37").Value = "9"
$ws.Range("G38").Value = "9"
$ws.Range("D39").Value = "0.01761"
$ws.Range("E39").Value = "-3.43%"
$ws.Range("G39").Value = "9"
$ws.Range("D40").Value = "0.04705"
$ws.Range("E40").Value = "4.19%"
$ws.Range("G40").Value = "9"
$ws.Range("D41").Value = "0.007202"
$ws.Range("E41").Value = "-0.56%"
$ws.Range("G41").Value = "9"
$ws.Range("D42").Value = "0.1362"
$ws.Range("E42").Value = "1.55%"
$ws.Range("G42").Value = "9"
$ws.Range("D43").Value = "0.002186"
$ws.Range("E43").Value = "-3.19%"
$ws.Range("G43").Value = "9"
$ws.Range("D44").Value = "0.009766"
$ws.Range("E44").Value = "-8.06%"
$ws.Range("G44").Value = "9"
$ws.Range("D45").Value = "0.00006252"
$ws.Range("E45").Value = "-0.79%"
$ws.Range("G45").Value = "9"
$ws.Range("E46").Value = "-0.69%"
$ws.Range("G46").Value = "9"
$ws.Range("E47").Value = "19.67%"
$ws.Range("G47").Value = "9"
$ws.Range("D48").Value = "0.7437"
$ws.Range("E48").Value = "-9.37%"
$ws.Range("G48").Value = "9"
$ws.Range("E49").Value = "-0.69%"
$ws.Range("G49").Value = "9"
$ws.Range("E50").Value = "-0.69%"
$ws.Range("G50").Value = "9"
$ws.Range("G51").Value = "9"

Write-Output "Updated crypto price/volume/hour data for rows 2-51"
